$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename sheet "My Series" -> "Data"
$ws.Name = "Data"

# 2. Update the number format used by A1 (and the new cells below) from
#    "0.000" to "###0.000", and extend it across the whole new row range.
$ws.Range("A1:J1").NumberFormat = "###0.000"

# 3. Populate the rest of row 1 (B1:J1) with the new data series values.
$ws.Range("B1").Value = 29205.22
$ws.Range("C1").Value = 22352.63
$ws.Range("D1").Value = 16421.36
$ws.Range("E1").Value = 10540.67
$ws.Range("F1").Value = 8354.02
$ws.Range("G1").Value = 7111.74
$ws.Range("H1").Value = 6384.48
$ws.Range("I1").Formula = "=NA()"
$ws.Range("J1").Value = 3271.34

# 4. Replace the (compressed/base64 encoded) CEIC add-in metadata payload
#    stored in the cell comment on A1 with the updated payload.
$c = $ws.Range("A1").Comment
[void]$c.Text("BxoAAB+LCAAAAAAAAAOlGdtuG8f1VxZ8aoGSu0vKtiSMN+BNClFSFEgqsvJSDHdH4lTLXWZ3VhTfUqBFijRFUQROkSZt0acUBeoabQKkdi//Eliy+9Rf6JnLzl5IxaVrCNbOuc2ZM+c6Qm9dz33jikQxDYOHFbtmVQwSuKFHg4uHlYSdV+37lbcc1L12iX+MIzwnDIgN4Ari/euYPqzMGFvsm+ZyuawtG7UwujDrlmWbjwb9sTsjc1ylQcxw4JKK5vJez1VxUNubDwjDHmZYcj6s9Ma9WptQtwOwAQ7wBYlqrSSmAYnjbsAooyTmnBHBjLQ7g3fkwZx67X7NRuYaPKNsJdT3JF2BUsIVHWxLJnROnLpl71at3WrDntj1fXtvf2e3tmfV300ZNSHq45iNSXRFXQEYMzxfCHZrt2HbdXtvZxeZG4lAVmYABw19b0SuaEy8NvH9eCuLmOoCmy6DU29nTAuZOV4l6M1VOIzwYjahzCfbqTEatIx5oHTJhDjoIIyIC/Z7I5WOyHIYKbNOFn3ATmY0YqsOXm0t6yQm0XDBjbQdq4M6YcCaPonYyQLumnjgCoBwWJQQZN6BzJg6NHbhmwYJ8Zxz7Md5pgISnYbRZbzALjmCODa5jGXgh9gDh2M0ZtSNMwFrGHQchQsQCbu3Qt87ALEp9QaMlt0LwMh841YYXpbFF5FI3Ku4YbjVOdby1+BoPAuXw8BfjZNp7EZ0SrxOK6XeiEM8JBV3O4lZOActMhCSsBxkBf8gBMtg1CEunWP/2Aczxk4DpBQAqJmw8Jyydugn80DbswRFp3CiCbnWJ9RrNITrDbjVw6AXlKRsxBU5RuEy9Zx1uDBCDtyMXX0j65gydQdgyifXEeJG+CkPqA8VIn8XOWjRK8YzQthGl5AYxJPhAa85TmvF90RmBkHgmuDfAHXsvb0HVcuGn4ll7Ysf2FmjUTfwxAdUmPtVu86zdkaXItFRMh9OIYavxJkcG3AlEIJT+C0fB5cAPaVsdtRMtd+AQfLMd9Kv4xBE68LHKwHWdsnDUC9w/cQjMgn0gnPhlFw3RX0nHq2B+hDXDsLBarJa8GRg3kHRIec48aE+MUgjF1mGKIFRM74s0+RB6CTy0zt2ePWPofy73rzmQoLkJa7mhnMOMKHqno6RmafnWd4l3eCij4OLBPKotmMZrj2MZ4BJhIOYH0enzZKzbSZC6b3IdO6sXy4ySyRoQuaLMML+AGxCD5JAVEtVDsC4A8xmagVx6xM3ta+ZsWquolKpzq8jExEoT8DTez4Cc1BBxR1eNhk5ogyI+CEHoQd1Hvt0GhU8bCMObiurfqmv8eNtWQnTC4BOE5LL98mKtx/ZQsG5wzp2ihAr1AaVnPFoZ7d+z2rUIXHzNRJHHhHsG13oQRkxesEVidkc2PaNEYmpB18U+/vG22RKKES8sJGKua2583zoICLvJdBLr4QqTdC3CCkSQDq9oAH21wk1JmNwzgiO/FWOUB61H7pAd/vTf918/vTFsy9uP3z86qsf/+fvv37xj1/cPPkJfNz+5a83H/1KHlMSowme+kQoNGnt7lqNHfA0DUIqN0CB9xKXCdjZmSj4eo1U5yoW7W6vfdhviWSigSk7byYSHg19vAqTbDmWhxAbiSs1U0+QJM4kTU5qXcB2CK/yMgx5k3pFitR5/F2M0hYvn3/x8vmf7uRWBivUnXu8nry+7tj2Gp2uO7z5l71dOlrsVK171Xo9R1yiQSOYkKCV03bqeU7DtmACaVi2TuSeduRNRGWUkjTBF2aJT4LaYRKwaKVdIL9OkcLxJxAiGi1DIbdQLvrVz179+XGBSllXQYpSQLkwiaQ3melCiD4aTYzx8GTU7hqT7pj7SYbL0Unh30KsdtfxVHCqIEiw/z2Yi/lUbFSg4asY4blBsDszVhCJuTgsONsmqNzoDUWWtTyMwmQhbyTHkEE3UOpsspFjQ64ROGHPtaSToTaQS11v/vblJgZ1EOVmJwFlesDLw1ABI0E5vIraT//54usPXjx7dvv0lzdf/6ggQe2j5x7wc4im/FK7PaQ8VW9KEHQ6Fsa8tH6Qqy8KyLvG45AGLHbs+6JhVCsErDaXJn6j3hxKnhAs7AXwEgS9jePuNVOB7RwhswgAPRcYqm2YtdYaIHN4Ztd//+a3t599efvJ01cf/PHmwz/cfPTJy+e/e/Xk9zLqbh8/vf35E5Xly4VA6MIbdtkBGmL+cw0ejQav3cY3739sBCEzoOkwEpGRvnn/05wwrqhoTzLJ0M9pRYoqrJHmmTmfkVNF61Dg0yyyAWjzEtbQFKqIhQvqZpu8W+WieNwJxHd6k2oSEyOEfuq7cJIiccb8v/IpFllSjx9YdbuusFIbfoQpjnOmP/TDKTQZKULMVyWSAte3M2S0Yr/D/rDV7GckUolh5ME4ZvEhkX+gtKnkJaUXpys9VWQQwELj5yY+H4nXyNZRWnIujZlqtDxvejz9OZYlhrdifSxQwHgeRbIhCtRj5ThZQDvM1ER6N168y+Qa4CPZq+Zb4mzd6xTxsM5hoRAW0Rwg8CI1KZRMU72Yj7GynT3ipsmWgCu85YA51Huk7LSuoK+MTJ53ulEURhuTT4ZJyQbQSUNGMTOLaxpxp7Lr9rK7SgFpwtMfcuxTJww7xCdsu8c6M+MehFdvzAt3vy1rLx76njLmdqOHNksmIP9iyR3l/32wlM7WjCJorPj7xtYvjOnUOoJhd0tt5FEEIx8BYXf1BHhAo5g94plAfUnImYacyQ71kbMn29BHci0I5Ic6oxJuFrRMI5fJZ+PQ79M53XIqtNLwLgoBUy4WsoPrbecovLIckWvoL3MSICdOfwhVg08820mT/gqpVPPzd5qYXszYtoo9mGLikalVdaekXt3xrN3qHiGNqm3D/9it1y3rHn/kUcIhcVCy3HITM72w7M85zn8BDlJ/rgcaAAA=")
